$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Aumento de limite en Descr Enfermedad y Descr Peste
$ws.Range("B12").Value = 350
$ws.Range("B22").Value = 350

# Update the active cell selection to D13
$ws.Range("D13").Select()
